$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the calendar dates to each week label in column A
$ws.Range("A2").Value = "Week 1: 10th Oct"
$ws.Range("A3").Value = "Week 2: 17th Oct"
$ws.Range("A4").Value = "Week 3: 24th Oct"
$ws.Range("A5").Value = "Week 4: 31st Oct"
$ws.Range("A6").Value = "Week 5: 7th Nov"
$ws.Range("A7").Value = "Week 6: 14th Nov"
$ws.Range("A8").Value = "Week 7: 21st Nov"
$ws.Range("A9").Value = "Week 8: 28th Nov"
$ws.Range("A10").Value = "Week 9: 5th Dec"
$ws.Range("A11").Value = "Week 10: 12th Dec"

# Widen column A to fit the longer labels
$ws.Columns.Item(1).ColumnWidth = 16.6

# Record Week 3 attendance (row 4) for the whole team, matching the
# formatting already used for the Week 1 / Week 2 rows above it
$ws.Range("B2:E2").Copy()
$ws.Range("B4:E4").PasteSpecial(-4122)
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
